$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) storage for the data range so that numeric-looking
# strings (e.g. "167.00", "1.00") are not silently coerced into numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '68.600.79'
$ws.Cells.Item(2, 5).Value = '  -0.56%  '
$ws.Cells.Item(3, 4).Value = '3.907.51'
$ws.Cells.Item(3, 5).Value = '  +2.70%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '601.89'
$ws.Cells.Item(5, 5).Value = '  +0.18%  '
$ws.Cells.Item(6, 4).Value = '167.00'
$ws.Cells.Item(6, 5).Value = '  +2.18%  '
$ws.Cells.Item(7, 4).Value = '3.908.38'
$ws.Cells.Item(7, 5).Value = '  +2.74%  '
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 5).Value = '  -0.95%  '
$ws.Cells.Item(10, 5).Value = '  -0.16%  '
$ws.Cells.Item(11, 5).Value = '  +2.54%  '
$ws.Cells.Item(12, 5).Value = '  +0.53%  '
$ws.Cells.Item(13, 4).Value = '0.0000257'
$ws.Cells.Item(13, 5).Value = '  +4.49%  '
$ws.Cells.Item(14, 4).Value = '37.37'
$ws.Cells.Item(14, 5).Value = '  +0.61%  '
$ws.Cells.Item(15, 4).Value = '4.558.56'
$ws.Cells.Item(15, 5).Value = '  +2.66%  '
$ws.Cells.Item(16, 4).Value = '3.889.66'
$ws.Cells.Item(16, 5).Value = '  +2.52%  '
$ws.Cells.Item(17, 4).Value = '68.703.84'
$ws.Cells.Item(17, 5).Value = '  -0.63%  '
$ws.Cells.Item(18, 4).Value = '7.47'
$ws.Cells.Item(18, 5).Value = '  +0.05%  '
$ws.Cells.Item(19, 4).Value = '17.36'
$ws.Cells.Item(19, 5).Value = '  +0.83%  '
$ws.Cells.Item(20, 4).Value = '0.112'
$ws.Cells.Item(20, 5).Value = '  -2.15%  '
$ws.Cells.Item(21, 4).Value = '11.06'
$ws.Cells.Item(21, 5).Value = '  -4.09%  '
$ws.Cells.Item(22, 4).Value = '492.00'
$ws.Cells.Item(22, 5).Value = '  +1.32%  '
$ws.Cells.Item(23, 4).Value = '0.728'
$ws.Cells.Item(23, 5).Value = '  +1.10%  '
$ws.Cells.Item(24, 4).Value = '0.0000166'
$ws.Cells.Item(24, 5).Value = '  +3.26%  '
$ws.Cells.Item(25, 4).Value = '84.80'
$ws.Cells.Item(25, 5).Value = '  +0.13%  '
$ws.Cells.Item(26, 4).Value = '2.23'
$ws.Cells.Item(26, 5).Value = '  -0.89%  '
$ws.Cells.Item(27, 4).Value = '12.05'
$ws.Cells.Item(27, 5).Value = '  -1.38%  '
$ws.Cells.Item(28, 4).Value = '10.16'
$ws.Cells.Item(28, 5).Value = '  +1.32%  '
$ws.Cells.Item(29, 5).Value = '  +0.16%  '
$ws.Cells.Item(30, 4).Value = '2.94'
$ws.Cells.Item(30, 5).Value = '  -0.70%  '
$ws.Cells.Item(31, 4).Value = '4.056.62'
$ws.Cells.Item(31, 5).Value = '  +2.60%  '
$ws.Cells.Item(32, 5).Value = '  -0.89%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = '31.86'
$ws.Cells.Item(33, 5).Value = '  +0.54%  '
$ws.Cells.Item(34, 2).Value = 'NEARProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(34, 4).Value = '7.71'
$ws.Cells.Item(34, 5).Value = '  -4.09%  '
$ws.Cells.Item(35, 4).Value = '3.861.90'
$ws.Cells.Item(35, 5).Value = '  +3.04%  '
$ws.Cells.Item(36, 5).Value = '  -0.22%  '
$ws.Cells.Item(37, 5).Value = '  +0.94%  '
$ws.Cells.Item(38, 5).Value = '  -0.24%  '
$ws.Cells.Item(39, 5).Value = '  +1.06%  '
$ws.Cells.Item(40, 4).Value = '3.22'
$ws.Cells.Item(40, 5).Value = '  +6.78%  '
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  -0.06%  '
$ws.Cells.Item(42, 5).Value = '  -0.65%  '
$ws.Cells.Item(43, 4).Value = '430.19'
$ws.Cells.Item(43, 5).Value = '  +0.21%  '
$ws.Cells.Item(44, 5).Value = '  +0.12%  '
$ws.Cells.Item(45, 4).Value = '48.05'
$ws.Cells.Item(45, 5).Value = '  -1.20%  '
$ws.Cells.Item(46, 4).Value = '8.57'
$ws.Cells.Item(46, 5).Value = '  +2.42%  '
$ws.Cells.Item(47, 5).Value = '  -0.01%  '
$ws.Cells.Item(48, 4).Value = '0.000275'
$ws.Cells.Item(48, 5).Value = '  +20.99%  '
$ws.Cells.Item(49, 4).Value = '142.88'
$ws.Cells.Item(49, 5).Value = '  +0.70%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '26.05'
$ws.Cells.Item(50, 5).Value = '  +4.59%  '
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '2.806.68'
$ws.Cells.Item(51, 5).Value = '  -0.73%  '

# Restore the default (unformatted) style now that the text values are locked in.
$dataRange.Style = "Normal"

